$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data to the right
$ws.Columns("A").Insert()

# Populate the new sNo column
$ws.Range("A1").Value = "sNo"
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2

# Match formatting used by the rest of the table (copy from neighboring column)
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("B2:B3").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A4").Select()
